$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "booked"
$ws.Range("C3").Value = 4.99
$ws.Range("D3").Value = 1

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "booked"
$ws.Range("C4").Value = 41.970000000000006
$ws.Range("D4").Value = 1

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "paid"
$ws.Range("C5").Value = 109.94999999999999
$ws.Range("D5").Value = 1
